$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colLetters = @('B','C','E','F','G','H','I','J','N','O')

$data = @(
    @(0.6896187935121247, 0.181058051359372, 0.824681771033525, 1.961027168877578, 0.2102521736780076, 0.3798410003148334, 0.2459880266992993, 0.02121218483690512, 0.8454711975907543, 1.090506345116282),
    @(0.6031622436052828, 0.1605233397285133, 0.7995253972065228, 1.925158862440043, 0.2086048924463739, 0.3825718963416307, 0.2505515470071007, 0.02126447437241907, 0.8373625661899098, 1.092489104206493),
    @(0.549900385240818, 0.147846999175556, 0.7844856960294919, 1.904360927396965, 0.2078194443181971, 0.3844689260907117, 0.2535616713022346, 0.02131490927722268, 0.8326978890959253, 1.094638321525338),
    @(0.5281528175523817, 0.1426645650038267, 0.7784590137853229, 1.896193213151548, 0.2075559685719455, 0.3852973329026312, 0.2548405271956327, 0.02134005753286949, 0.8308763912083492, 1.095747901541657),
    @(0.5245391057170252, 0.1418030264514698, 0.7774644532850346, 1.894855534682762, 0.2075156305410104, 0.3854382314176519, 0.2550560297153677, 0.02134451055539621, 0.8305787407633574, 1.095946248634419),
    @(0.5496072615924277, 0.1477771742834477, 0.784404004715995, 1.90424952991313, 0.2078156621183425, 0.3844798742000037, 0.2535787071852162, 0.02131522984560164, 0.8326730017018775, 1.094652339971262),
    @(0.6598462993849523, 0.1739919960422469, 0.8159234448559971, 1.948405130557305, 0.2096371499643297, 0.380736881404502, 0.247518262540968, 0.02122640202388482, 0.842610407761029, 1.090996271924283),
    @(0.8745574503307694, 0.2248464778054995, 0.880965424179152, 2.044747073888118, 0.2150129461363903, 0.3751460650859855, 0.2372903815011469, 0.02119827194665902, 0.8645728497982219, 1.091246509003398),
    @(1.031343672512889, 0.2618575020689775, 0.9307391682171726, 2.121527673911743, 0.2200778733138975, 0.3721073359751728, 0.2307931460565342, 0.02126755848671991, 0.882198092454999, 1.095992412247028),
    @(1.102447494855369, 0.2786154676547312, 0.953818108945967, 2.157771605054961, 0.2226276214609868, 0.3709575909952889, 0.2280598073903679, 0.02131880979904466, 0.8905358197910687, 1.099150512718239),
    @(1.129339713133277, 0.2849496570157442, 0.9626204889777199, 2.171686243616676, 0.2236287309327452, 0.3705557018316199, 0.2270568606606158, 0.0213410701928396, 0.8937387369865633, 1.100490722879698),
    @(1.123549499831256, 0.2835860021330063, 0.9607219404231557, 2.16868102405283, 0.2234115379790609, 0.3706407654132704, 0.2272714328112979, 0.02133614890645319, 0.8930469097048501, 1.100195656623043),
    @(1.104660610412054, 0.2791368217408206, 0.9545410251277104, 2.158912562438189, 0.222709268956649, 0.3709238557789405, 0.2279766500126463, 0.02132058391726943, 0.8907984140949026, 1.099257876920689),
    @(1.093086237117461, 0.2764100363827708, 0.9507632287176619, 2.152953839435185, 0.2222837491152774, 0.37110162018719, 0.2284128017829854, 0.0213114218579662, 0.8894270723254323, 1.098702271266291),
    @(1.026692312603473, 0.2607607154647553, 0.9292396995713119, 2.119185581455923, 0.2199162061549345, 0.3721871589180523, 0.2309762615196824, 0.02126460745450132, 0.8816596053362389, 1.095806182382972),
    @(0.9859044397693992, 0.2511399676254484, 0.9161476103608948, 2.09880734834266, 0.2185268967590446, 0.3729127075938408, 0.2326058970113607, 0.02124095306884755, 0.8769761657304684, 1.094285882827648),
    @(0.9624238242038246, 0.2455989993918024, 0.9086584769450354, 2.08721016438102, 0.2177509133753688, 0.3733519133951972, 0.2335641342438954, 0.02122920410718621, 0.8743125108588004, 1.093505455514133),
    @(0.9544702224038701, 0.2437216692676714, 0.9061298392627464, 2.08330480796451, 0.217492140508277, 0.3735043789848049, 0.2338921644900402, 0.02122554449079317, 0.8734158312923483, 1.093257343192221),
    @(0.9902485139386386, 0.2521648781394674, 0.9175370322184051, 2.100963825054663, 0.218672397469021, 0.3728332058774555, 0.2324302539605689, 0.02124327887719346, 0.8774716098726287, 1.094437987201928),
    @(1.110209651682055, 0.2804439742391764, 0.9563548022458406, 2.161776639158404, 0.2229145749386419, 0.3708397958637875, 0.2277686380359949, 0.02132507819203511, 0.8914576173483795, 1.099529404137002),
    @(1.188416534225496, 0.2988576861714876, 0.9820910561452934, 2.202628192550122, 0.2258945558369732, 0.3697322401688297, 0.2249092271373563, 0.02139517311530525, 0.9008638952505805, 1.103698355240937),
    @(1.146694440879799, 0.2890363224082648, 0.9683215551588233, 2.180723471743818, 0.2242850192809129, 0.3703054818090266, 0.2264181697141439, 0.02135623509837359, 0.895819415650223, 1.101396111226478),
    @(0.9882846548190969, 0.2517015468175146, 0.9169087572242489, 2.099988512747899, 0.2186065457542838, 0.3728690798336771, 0.2325095957275263, 0.02124222161729605, 0.8772475294699262, 1.094368929259929),
    @(0.8166364822563992, 0.2111497004351861, 0.824681771033525, 2.017634490371265, 0.2133639143985064, 0.3764710097767718, 0.239879189275781, 0.02119015085732912, 0.8583684792001662, 1.090380463252529)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "{0}{1}" -f $colLetters[$j], $rowNum
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
